# Update cryptocurrency market data (prices & 1h volume %) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''68.818.04'
$ws.Range("E2").Value = '  +0.66%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '''3.848.74'
$ws.Range("E3").Value = '  -1.44%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.05%  '

# Row 5: BNB
$ws.Range("D5").Value = '''523.47'
$ws.Range("E5").Value = '  +7.44%  '

# Row 6: Solana
$ws.Range("D6").Value = '''142.22'
$ws.Range("E6").Value = '  -2.83%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.606'
$ws.Range("E7").Value = '  -2.74%  '

# Row 8: USDC
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  +0.10%  '

# Row 9: Cardano
$ws.Range("D9").Value = '''0.711'
$ws.Range("E9").Value = '  -4.49%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  -6.70%  '

# Row 11: ShibaInu
$ws.Range("E11").Value = '  -8.13%  '

# Row 12: Avalanche
$ws.Range("D12").Value = '''41.60'
$ws.Range("E12").Value = '  -3.81%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '''4.464.90'
$ws.Range("E13").Value = '  -1.19%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''10.12'
$ws.Range("E14").Value = '  -3.39%  '

# Row 15: WrappedEther
$ws.Range("D15").Value = '''3.888.33'
$ws.Range("E15").Value = '  -0.12%  '

# Row 16: Uniswap
$ws.Range("D16").Value = '''13.85'
$ws.Range("E16").Value = '  -2.67%  '

# Row 17: 'TRON' -> 'Chainlink'
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = '''20.44'
$ws.Range("E17").Value = '  +2.08%  '

# Row 18: 'Polygon' -> 'TRON'
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '''0.134'
$ws.Range("E18").Value = '  -1.55%  '

# Row 19: 'Chainlink' -> 'Polygon'
$ws.Range("B19").Value = 'Polygon'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D19").Value = '''1.21'
$ws.Range("E19").Value = '  +4.81%  '

# Row 20: WrappedBTC
$ws.Range("D20").Value = '''68.769.98'
$ws.Range("E20").Value = '  +0.49%  '

# Row 21: BitcoinCash
$ws.Range("D21").Value = '''420.00'
$ws.Range("E21").Value = '  -2.95%  '

# Row 22: ImmutableX
$ws.Range("E22").Value = '  -5.10%  '

# Row 23: InternetComputer(DFINITY)
$ws.Range("E23").Value = '  -4.63%  '

# Row 24: Litecoin
$ws.Range("D24").Value = '''86.85'
$ws.Range("E24").Value = '  -3.25%  '

# Row 25: PancakeSwap
$ws.Range("E25").Value = '  +5.45%  '

# Row 26: RenderToken
$ws.Range("D26").Value = '''11.29'
$ws.Range("E26").Value = '  -7.74%  '

# Row 27: Filecoin
$ws.Range("D27").Value = '''10.50'
$ws.Range("E27").Value = '  -4.37%  '

# Row 28: EthereumClassic
$ws.Range("D28").Value = '''35.92'
$ws.Range("E28").Value = '  -4.37%  '

# Row 29: Bittensor
$ws.Range("D29").Value = '''684.04'
$ws.Range("E29").Value = '  -4.05%  '

# Row 30: Cosmos
$ws.Range("E30").Value = '  -2.62%  '

# Row 31: Hedera
$ws.Range("E31").Value = '  -5.09%  '

# Row 32: Toncoin
$ws.Range("E32").Value = '  -3.31%  '

# Row 33: OKB
$ws.Range("D33").Value = '''67.98'
$ws.Range("E33").Value = '  +10.46%  '

# Row 34: TheGraph
$ws.Range("E34").Value = '  +6.54%  '

# Row 35: NEARProtocol
$ws.Range("D35").Value = '''5.88'
$ws.Range("E35").Value = '  -3.29%  '

# Row 36: PEPE
$ws.Range("D36").Value = '''0.0₃0848'
$ws.Range("E36").Value = '  -5.90%  '

# Row 37: InjectiveProtocol
$ws.Range("D37").Value = '''39.62'
$ws.Range("E37").Value = '  -3.06%  '

# Row 38: Dai
$ws.Range("E38").Value = '  +0.08%  '

# Row 39: Kaspa
$ws.Range("E39").Value = '  -0.73%  '

# Row 40: FirstDigitalUSD
$ws.Range("E40").Value = '  -0.18%  '

# Row 41: ThetaToken
$ws.Range("D41").Value = '''3.21'
$ws.Range("E41").Value = '  +3.87%  '

# Row 42: VeChain
$ws.Range("D42").Value = '''0.0478'
$ws.Range("E42").Value = '  -3.68%  '

# Row 43: WEMIXToken
$ws.Range("D43").Value = '''3.14'
$ws.Range("E43").Value = '  +3.02%  '

# Row 44: Fetch.AI
$ws.Range("D44").Value = '''2.75'
$ws.Range("E44").Value = '  -6.95%  '

# Row 45: ApeXProtocol
$ws.Range("D45").Value = '''3.40'
$ws.Range("E45").Value = '  +0.79%  '

# Row 46: Stellar
$ws.Range("E46").Value = '  -2.93%  '

# Row 47: Stacks
$ws.Range("D47").Value = '''2.94'
$ws.Range("E47").Value = '  +4.60%  '

# Row 48: Maker
$ws.Range("D48").Value = '''2.751.45'
$ws.Range("E48").Value = '  +14.22%  '

# Row 49: 'Monero' -> 'FLOKI'
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").Value = '''0.000270'
$ws.Range("E49").Value = '  +10.07%  '

# Row 50: 'FLOKI' -> 'Monero'
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '''144.22'
$ws.Range("E50").Value = '  +0.99%  '

# Row 51: 'BabyDogeCoin' -> 'LidoDAOToken'
$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").Value = '''3.25'
$ws.Range("E51").Value = '  -3.76%  '

